# E4_Fiche_Mathieu_Soufflard.docx edit script
# 1) "Descriptif détailler" -> "Descriptif détaillé" (typo fix)
# 2) Word re-materialised the (previously latent) built-in "heading 1" /
#    "heading 3" paragraph styles plus their linked "... Car" character
#    styles into styles.xml.

$d = $word.ActiveDocument

# --- 1. text fix -----------------------------------------------------
$d.Content.Find.Execute("détailler", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "détaillé", 2) | Out-Null

# --- 2. add the "Titre 1" / "heading 1" paragraph style ---------------
$titre1 = $d.Styles.Add("Titre1", 1)
$titre1.NameLocal = "heading 1"
$titre1.BaseStyle = "Normal"
$titre1.NextParagraphStyle = "Normal"
$titre1.Priority = 9
$titre1.QuickStyle = $true
$titre1.ParagraphFormat.KeepWithNext = $true
$titre1.ParagraphFormat.KeepTogether = $true
$titre1.ParagraphFormat.SpaceBefore = 12
$titre1.ParagraphFormat.OutlineLevel = 1
$titre1.Font.Size = 16
$titre1.Font.SizeBi = 16
$titre1.Font.TextColor.ObjectThemeColor = 4

# --- 3. add the "Titre 3" / "heading 3" paragraph style ---------------
$titre3 = $d.Styles.Add("Titre3", 1)
$titre3.NameLocal = "heading 3"
$titre3.BaseStyle = "Normal"
$titre3.NextParagraphStyle = "Normal"
$titre3.Priority = 9
$titre3.UnhideWhenUsed = $true
$titre3.QuickStyle = $true
$titre3.ParagraphFormat.KeepWithNext = $true
$titre3.ParagraphFormat.KeepTogether = $true
$titre3.ParagraphFormat.SpaceBefore = 2
$titre3.ParagraphFormat.OutlineLevel = 3
$titre3.Font.SizeBi = 12
$titre3.Font.TextColor.ObjectThemeColor = 4

# --- 4. linked character style for Titre 1 -----------------------------
$titre1Car = $d.Styles.Add("Titre1Car", 2)
$titre1Car.NameLocal = "Titre 1 Car"
$titre1Car.BaseStyle = "Policepardfaut"
$titre1Car.Priority = 9
$titre1Car.Font.Size = 16
$titre1Car.Font.SizeBi = 16
$titre1Car.Font.TextColor.ObjectThemeColor = 4

$titre1.LinkStyle = "Titre1Car"
$titre1Car.LinkStyle = "Titre1"

# --- 5. linked character style for Titre 3 -----------------------------
$titre3Car = $d.Styles.Add("Titre3Car", 2)
$titre3Car.NameLocal = "Titre 3 Car"
$titre3Car.BaseStyle = "Policepardfaut"
$titre3Car.Priority = 9
$titre3Car.Font.Size = 12
$titre3Car.Font.SizeBi = 12
$titre3Car.Font.TextColor.ObjectThemeColor = 4

$titre3.LinkStyle = "Titre3Car"
$titre3Car.LinkStyle = "Titre3"

Write-Host "done"
